# Update "want to go" counts (column F) that were refreshed in the source data.
# Affects both the "展览" sheet and the "全部类型" aggregate sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 568
    $ws.Range("F8").Value = 485
    $ws.Range("F9").Value = 3598
}
